# "Add files via upload" — the re-uploaded workbook keeps the same three
# lookup tables (Links / Deel / Rechts) but with their helper "year"
# formulas (column B, rows 2-4) stripped out, leaving only the shared
# title text in column A (and the column B header in row 1). The active
# sheet/selection also moved: "Links" is now the selected tab (cell B13),
# while "Deel" and "Rechts" are left with B2:B4 selected.

$wb = $excel.ActiveWorkbook

$wsLinks  = $wb.Worksheets.Item("Links")
$wsDeel   = $wb.Worksheets.Item("Deel")
$wsRechts = $wb.Worksheets.Item("Rechts")

# Remove the helper formulas from column B (rows 2-4) on every sheet.
$null = $wsLinks.Range("B2:B4").ClearContents()
$null = $wsDeel.Range("B2:B4").ClearContents()
$null = $wsRechts.Range("B2:B4").ClearContents()

# Restore each sheet's selection, then make "Links" the active/selected tab.
$null = $wsDeel.Range("B2:B4").Select()
$null = $wsRechts.Range("B2:B4").Select()

$null = $wsLinks.Activate()
$null = $wsLinks.Range("B13").Select()
